# Generate Report for Handoff
# The file "8c9ec3f7-c2f4-4aca-8b20-1c5c32e828a9.md" moved from
# "Handed back: in sync with en-US" to "Ready for handoff" with new
# handoff timestamps, for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet: row 3 is the 8c9ec3f7-... file ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = "2016-03-24 19:02:55"

# --- zh-cn sheet: row 3 is the 8c9ec3f7-... file ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("E3").Value = "2016-03-24 19:02:51"

# --- de-de sheet: row 3 is the 8c9ec3f7-... file ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("E3").Value = "2016-03-24 19:02:55"
